# Weekly update: insert a new price record for "Poroto granado" at the top
# of the data block (row 62), pushing every subsequent row down by one.
# The previously-last row (old row 192) lands on the new row 193, and the
# sheet's used range grows from A1:R192 to A1:R193.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record.
$ws.Rows(62).Insert()

# Populate the newly inserted row with the latest week's data.
$ws.Range("A62").Value = 5
$ws.Range("B62").Value = "Macroferia Regional de Talca"
$ws.Range("C62").Value = "Maule"
$ws.Range("D62").Value = 44979
$ws.Range("E62").Value = 7
$ws.Range("F62").Value = 100112030
$ws.Range("G62").Value = "Poroto granado"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 200
$ws.Range("K62").Value = 25000
$ws.Range("L62").Value = 25000
$ws.Range("M62").Value = 25000
$ws.Range("N62").Value = "$/saco 25 kilos"
$ws.Range("O62").Value = "Región del Maule"
$ws.Range("P62").Value = 1000
$ws.Range("Q62").Value = 25
$ws.Range("R62").Value = "Hortaliza"
